$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$d.Paragraphs.Item(1).Range.Find.Execute("2023-08-07 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-08 Tuesday", 2) | Out-Null

# Update each arithmetic-problem cell in the table, addressed by (row, column)
# to avoid any ambiguity from duplicate cell text (e.g. "39+59=" appears twice).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "31+57="
$t.Cell(1, 2).Range.Text = "1+33="
$t.Cell(1, 3).Range.Text = "82-18="
$t.Cell(1, 4).Range.Text = "89+5="
$t.Cell(1, 5).Range.Text = "67-42="

$t.Cell(2, 1).Range.Text = "35+5="
$t.Cell(2, 2).Range.Text = "35+11="
$t.Cell(2, 3).Range.Text = "45-18="
$t.Cell(2, 4).Range.Text = "90-12="
$t.Cell(2, 5).Range.Text = "9+4="

$t.Cell(3, 1).Range.Text = "96-94="
$t.Cell(3, 2).Range.Text = "87-66="
$t.Cell(3, 3).Range.Text = "57+10="
$t.Cell(3, 4).Range.Text = "73+2="
$t.Cell(3, 5).Range.Text = "39+10="

$t.Cell(4, 1).Range.Text = "27-7="
$t.Cell(4, 2).Range.Text = "97-51="
$t.Cell(4, 3).Range.Text = "16+30="
$t.Cell(4, 4).Range.Text = "76-18="
$t.Cell(4, 5).Range.Text = "13-7="

$t.Cell(5, 1).Range.Text = "83-19="
$t.Cell(5, 2).Range.Text = "35-12="
$t.Cell(5, 3).Range.Text = "24-16="
$t.Cell(5, 4).Range.Text = "55-8="
$t.Cell(5, 5).Range.Text = "7+14="

$t.Cell(6, 1).Range.Text = "26+43="
$t.Cell(6, 2).Range.Text = "98-55="
$t.Cell(6, 3).Range.Text = "95-0="
$t.Cell(6, 4).Range.Text = "77-31="
$t.Cell(6, 5).Range.Text = "54-6="

$t.Cell(7, 1).Range.Text = "13-5="
$t.Cell(7, 2).Range.Text = "70-4="
$t.Cell(7, 3).Range.Text = "30+65="
$t.Cell(7, 4).Range.Text = "27+37="
$t.Cell(7, 5).Range.Text = "34-7="

$t.Cell(8, 1).Range.Text = "30+32="
$t.Cell(8, 2).Range.Text = "55-15="
$t.Cell(8, 3).Range.Text = "32+50="
$t.Cell(8, 4).Range.Text = "32+59="
$t.Cell(8, 5).Range.Text = "1+46="

$t.Cell(9, 1).Range.Text = "45-1="
$t.Cell(9, 2).Range.Text = "31+22="
$t.Cell(9, 3).Range.Text = "78-29="
$t.Cell(9, 4).Range.Text = "16-13="
$t.Cell(9, 5).Range.Text = "31-9="

$t.Cell(10, 1).Range.Text = "56+28="
$t.Cell(10, 2).Range.Text = "33+42="
$t.Cell(10, 3).Range.Text = "77+13="
$t.Cell(10, 4).Range.Text = "88-32="
$t.Cell(10, 5).Range.Text = "20+57="

$t.Cell(11, 1).Range.Text = "0+66="
$t.Cell(11, 2).Range.Text = "80-22="
$t.Cell(11, 3).Range.Text = "3+38="
$t.Cell(11, 4).Range.Text = "71-18="
$t.Cell(11, 5).Range.Text = "44-16="

$t.Cell(12, 1).Range.Text = "23+53="
$t.Cell(12, 2).Range.Text = "83-47="
$t.Cell(12, 3).Range.Text = "84-14="
$t.Cell(12, 4).Range.Text = "23+1="
$t.Cell(12, 5).Range.Text = "64-40="

$t.Cell(13, 1).Range.Text = "49-16="
$t.Cell(13, 2).Range.Text = "85-13="
$t.Cell(13, 3).Range.Text = "63-1="
$t.Cell(13, 4).Range.Text = "4+90="
$t.Cell(13, 5).Range.Text = "84-48="

$t.Cell(14, 1).Range.Text = "71+5="
$t.Cell(14, 2).Range.Text = "63+16="
$t.Cell(14, 3).Range.Text = "13-5="
$t.Cell(14, 4).Range.Text = "7+20="
$t.Cell(14, 5).Range.Text = "48-28="

$t.Cell(15, 1).Range.Text = "95-18="
$t.Cell(15, 2).Range.Text = "34-30="
$t.Cell(15, 3).Range.Text = "21+21="
$t.Cell(15, 4).Range.Text = "80-30="
$t.Cell(15, 5).Range.Text = "65+5="

$t.Cell(16, 1).Range.Text = "3+64="
$t.Cell(16, 2).Range.Text = "30+2="
$t.Cell(16, 3).Range.Text = "84+15="
$t.Cell(16, 4).Range.Text = "2+27="
$t.Cell(16, 5).Range.Text = "19+47="

$t.Cell(17, 1).Range.Text = "99-46="
$t.Cell(17, 2).Range.Text = "22-17="
$t.Cell(17, 3).Range.Text = "57-7="
$t.Cell(17, 4).Range.Text = "30-18="
$t.Cell(17, 5).Range.Text = "82+5="

$t.Cell(18, 1).Range.Text = "60-43="
$t.Cell(18, 2).Range.Text = "14+33="
$t.Cell(18, 3).Range.Text = "83-6="
$t.Cell(18, 4).Range.Text = "87-77="
$t.Cell(18, 5).Range.Text = "44+46="

$t.Cell(19, 1).Range.Text = "71-17="
$t.Cell(19, 2).Range.Text = "54-11="
$t.Cell(19, 3).Range.Text = "60-4="
$t.Cell(19, 4).Range.Text = "92-58="
$t.Cell(19, 5).Range.Text = "86-2="

$t.Cell(20, 1).Range.Text = "49+17="
$t.Cell(20, 2).Range.Text = "64+32="
$t.Cell(20, 3).Range.Text = "3+40="
$t.Cell(20, 4).Range.Text = "36-5="
$t.Cell(20, 5).Range.Text = "95-83="
